# Update sample data for session 11/12/2021:
# Column A (test_number) for rows 2-11 changes from 3 to 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A11").Value = 2

# Move the active selection to C14 (as last left by the editor).
$ws.Range("C14").Select()
